$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update Image Height (C3) and Image Width (D3): 227 -> 224
$ws.Range("C3").Value = 224
$ws.Range("D3").Value = 224

# Update Zero Padding (A11): 2 -> 5
$ws.Range("A11").Value = 5

# Update the active selection to E10
$ws.Range("E10").Select()
